$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as plain TEXT (matches the source data, which is
# stored as inline/shared strings), avoiding Excel's automatic number/date
# coercion for numeric-looking strings (e.g. "608.30", "0.999", "3.90").
# The cell's original style is restored afterward so no stray "Text"
# number-format style lingers on the cell.
function Set-TextValue {
    param(
        [string]$Address,
        [string]$Text
    )
    $cell = $ws.Range($Address)
    $originalStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $Text
    $cell.Style = $originalStyle
}

Set-TextValue "D2" '68.328.34'
Set-TextValue "E2" '  -0.09%  '
Set-TextValue "D3" '2.707.56'
Set-TextValue "E3" '  +2.29%  '
Set-TextValue "E4" '  +0.02%  '
Set-TextValue "D5" '608.30'
Set-TextValue "E5" '  +2.01%  '
Set-TextValue "D6" '166.39'
Set-TextValue "E6" '  +4.82%  '
Set-TextValue "E7" '  +0.02%  '
Set-TextValue "D8" '0.558'
Set-TextValue "E8" '  +3.44%  '
Set-TextValue "D9" '2.706.73'
Set-TextValue "E9" '  +2.28%  '
Set-TextValue "E10" '  +1.75%  '
Set-TextValue "E12" '  +3.32%  '
Set-TextValue "E13" '  +0.46%  '
Set-TextValue "D14" '28.38'
Set-TextValue "E14" '  +1.34%  '
Set-TextValue "D15" '3.200.91'
Set-TextValue "E15" '  +2.25%  '
Set-TextValue "E16" '  -0.06%  '
Set-TextValue "D17" '68.248.67'
Set-TextValue "E17" '  -0.07%  '
Set-TextValue "D18" '2.709.61'
Set-TextValue "E18" '  +2.51%  '
Set-TextValue "D19" '11.79'
Set-TextValue "E19" '  +2.14%  '
Set-TextValue "D20" '369.79'
Set-TextValue "E20" '  +1.76%  '
Set-TextValue "D21" '7.61'
Set-TextValue "E21" '  +1.94%  '
Set-TextValue "D22" '4.48'
Set-TextValue "E22" '  +1.65%  '
Set-TextValue "E23" '  +3.51%  '
Set-TextValue "E24" '  -0.61%  '
Set-TextValue "D25" '73.04'
Set-TextValue "E25" '  -2.29%  '
Set-TextValue "E26" '  +0.09%  '
Set-TextValue "D27" '9.98'
Set-TextValue "E27" '  +0.92%  '
Set-TextValue "E28" '  +2.29%  '
Set-TextValue "E29" '  +0.85%  '
Set-TextValue "E30" '  +0.20%  '
Set-TextValue "D31" '578.69'
Set-TextValue "E31" '  +2.33%  '
Set-TextValue "D32" '8.11'
Set-TextValue "E32" '  +0.76%  '
Set-TextValue "D33" '1.42'
Set-TextValue "E33" '  +1.31%  '
Set-TextValue "E34" '  +5.58%  '
Set-TextValue "E35" '  +1.75%  '
Set-TextValue "D36" '0.999'
Set-TextValue "E36" '  -0.03%  '
Set-TextValue "E37" '  -3.44%  '
Set-TextValue "D38" '162.50'
Set-TextValue "E38" '  +1.06%  '
Set-TextValue "E39" '  +0.84%  '
Set-TextValue "E40" '  +1.73%  '
Set-TextValue "B41" 'RenderToken'
Set-TextValue "C41" 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
Set-TextValue "D41" '5.37'
Set-TextValue "E41" '  +1.21%  '
Set-TextValue "B42" 'Stacks'
Set-TextValue "C42" 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
Set-TextValue "D42" '1.86'
Set-TextValue "E42" '  -0.35%  '
Set-TextValue "D43" '17.99'
Set-TextValue "E43" '  +0.98%  '
Set-TextValue "E44" '  -1.52%  '
Set-TextValue "E45" '  -0.04%  '
Set-TextValue "D46" '0.0₆0308'
Set-TextValue "E46" '  -3.78%  '
Set-TextValue "D47" '40.74'
Set-TextValue "E47" '  +1.06%  '
Set-TextValue "D48" '0.595'
Set-TextValue "E48" '  +3.64%  '
Set-TextValue "D49" '154.30'
Set-TextValue "E49" '  -2.46%  '
Set-TextValue "D50" '3.90'
Set-TextValue "E50" '  +2.09%  '
Set-TextValue "D51" '1.77'
Set-TextValue "E51" '  +4.10%  '
